# Auto-update data + news
# Update row 9 (ICSA_thou / Initial Jobless Claims) with refreshed series values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E9").Value = 209000
$ws.Range("G9").Value = 364051.724137931
$ws.Range("H9").Value = -13000
$ws.Range("I9").Value = -0.05855855855855856
